$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card12")

# Update header M1: remove trailing space -> "Serviced by"
$ws.Cells.Item(1, 13).Value = "Serviced by"

# Add new header N1: "Event " (with trailing space), matching style of M1
# (bold font, thin border, centered/top-aligned - same visual style as the
# other header cells in row 1)
$ws.Cells.Item(1, 14).Value = "Event "
$ws.Cells.Item(1, 14).Font.Bold = $true
$ws.Cells.Item(1, 14).HorizontalAlignment = -4108
$ws.Cells.Item(1, 14).VerticalAlignment = -4160
$ws.Cells.Item(1, 14).Borders.LineStyle = 1

# Fill M2:M13 with "nan" (these cells were present but empty)
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"
}

# Add new N2:N13 cells as empty (present, no value), matching style of column M data cells
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 14).Style = $ws.Cells.Item($r, 13).Style
}
